# Apply updated market-price / profit figures for several Leve rows
# across multiple crafting-profession sheets (scheduled market refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 783535
$ws.Cells.Item(15, 9).Value = 783535
$ws.Cells.Item(15, 11).Value = 2350605
$ws.Cells.Item(15, 13).Value = -2350436
$ws.Cells.Item(43, 8).Value = 7456.522
$ws.Cells.Item(43, 9).Value = 7647.0586
$ws.Cells.Item(43, 10).Value = 6916.6665
$ws.Cells.Item(43, 11).Value = 7647.0586
$ws.Cells.Item(43, 12).Value = 6916.6665
$ws.Cells.Item(43, 13).Value = -7578.0586
$ws.Cells.Item(43, 14).Value = -7054.6665
$ws.Cells.Item(58, 8).Value = 4706.636
$ws.Cells.Item(58, 9).Value = 789
$ws.Cells.Item(58, 10).Value = 11562.5
$ws.Cells.Item(58, 11).Value = 2367
$ws.Cells.Item(58, 12).Value = 34687.5
$ws.Cells.Item(58, 13).Value = -2217
$ws.Cells.Item(58, 14).Value = -34987.5
$ws.Cells.Item(62, 8).Value = 2500
$ws.Cells.Item(62, 9).Value = 2000
$ws.Cells.Item(62, 11).Value = 2000
$ws.Cells.Item(62, 13).Value = -1376
$ws.Cells.Item(65, 8).Value = 2500
$ws.Cells.Item(65, 9).Value = 2000
$ws.Cells.Item(65, 11).Value = 10000
$ws.Cells.Item(65, 13).Value = -6880
$ws.Cells.Item(112, 8).Value = 2699.6667
$ws.Cells.Item(112, 9).Value = 2100
$ws.Cells.Item(112, 10).Value = 2999.5
$ws.Cells.Item(112, 11).Value = 6300
$ws.Cells.Item(112, 12).Value = 8998.5
$ws.Cells.Item(112, 13).Value = -5192
$ws.Cells.Item(112, 14).Value = -11214.5
$ws.Cells.Item(132, 8).Value = 2075.25
$ws.Cells.Item(132, 9).Value = 2077.6775
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 6233.032499999999
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -3703.032499999999
$ws.Cells.Item(132, 14).Value = -11060

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 589.4211
$ws.Cells.Item(5, 9).Value = 364.70587
$ws.Cells.Item(5, 10).Value = 2499.5
$ws.Cells.Item(5, 11).Value = 364.70587
$ws.Cells.Item(5, 12).Value = 2499.5
$ws.Cells.Item(5, 13).Value = -252.70587
$ws.Cells.Item(5, 14).Value = -2723.5
$ws.Cells.Item(32, 8).Value = 3065.4521
$ws.Cells.Item(32, 9).Value = 2536.8208
$ws.Cells.Item(32, 11).Value = 2536.8208
$ws.Cells.Item(32, 13).Value = -2249.8208
$ws.Cells.Item(88, 8).Value = 1568.5
$ws.Cells.Item(88, 9).Value = 1412
$ws.Cells.Item(88, 10).Value = 1725
$ws.Cells.Item(88, 11).Value = 1412
$ws.Cells.Item(88, 12).Value = 1725
$ws.Cells.Item(88, 13).Value = -1006
$ws.Cells.Item(88, 14).Value = -2537
$ws.Cells.Item(91, 8).Value = 1568.5
$ws.Cells.Item(91, 9).Value = 1412
$ws.Cells.Item(91, 10).Value = 1725
$ws.Cells.Item(91, 11).Value = 1412
$ws.Cells.Item(91, 12).Value = 1725
$ws.Cells.Item(91, 13).Value = -8
$ws.Cells.Item(91, 14).Value = -4533

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 589.4211
$ws.Cells.Item(4, 9).Value = 364.70587
$ws.Cells.Item(4, 10).Value = 2499.5
$ws.Cells.Item(4, 11).Value = 364.70587
$ws.Cells.Item(4, 12).Value = 2499.5
$ws.Cells.Item(4, 13).Value = -249.70587
$ws.Cells.Item(4, 14).Value = -2729.5
$ws.Cells.Item(94, 8).Value = 722.65
$ws.Cells.Item(94, 9).Value = 627.6667
$ws.Cells.Item(94, 11).Value = 627.6667
$ws.Cells.Item(94, 13).Value = -176.6667
$ws.Cells.Item(115, 8).Value = 199950
$ws.Cells.Item(115, 10).Value = 199950
$ws.Cells.Item(115, 12).Value = 199950
$ws.Cells.Item(115, 14).Value = -203084

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 5205.15
$ws.Cells.Item(62, 9).Value = 4673.3335
$ws.Cells.Item(62, 10).Value = 6800.6
$ws.Cells.Item(62, 11).Value = 4673.3335
$ws.Cells.Item(62, 12).Value = 6800.6
$ws.Cells.Item(62, 13).Value = -4049.3335
$ws.Cells.Item(62, 14).Value = -8048.6
$ws.Cells.Item(65, 8).Value = 5205.15
$ws.Cells.Item(65, 9).Value = 4673.3335
$ws.Cells.Item(65, 10).Value = 6800.6
$ws.Cells.Item(65, 11).Value = 23366.6675
$ws.Cells.Item(65, 12).Value = 34003
$ws.Cells.Item(65, 13).Value = -20246.6675
$ws.Cells.Item(65, 14).Value = -40243
$ws.Cells.Item(105, 8).Value = 2458.6875
$ws.Cells.Item(105, 9).Value = 2146.4
$ws.Cells.Item(105, 10).Value = 2979.1667
$ws.Cells.Item(105, 11).Value = 2146.4
$ws.Cells.Item(105, 12).Value = 2979.1667
$ws.Cells.Item(105, 13).Value = -399.4000000000001
$ws.Cells.Item(105, 14).Value = -6473.1667
$ws.Cells.Item(132, 8).Value = 2917.1025
$ws.Cells.Item(132, 9).Value = 2307.6858
$ws.Cells.Item(132, 11).Value = 6923.057400000001
$ws.Cells.Item(132, 13).Value = -4393.057400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 423.75
$ws.Cells.Item(8, 9).Value = 423.75
$ws.Cells.Item(8, 11).Value = 1271.25
$ws.Cells.Item(8, 13).Value = -1132.25
$ws.Cells.Item(98, 8).Value = 3024.3333
$ws.Cells.Item(98, 9).Value = 2241.1667
$ws.Cells.Item(98, 10).Value = 3807.5
$ws.Cells.Item(98, 11).Value = 6723.500100000001
$ws.Cells.Item(98, 12).Value = 11422.5
$ws.Cells.Item(98, 13).Value = -5225.500100000001
$ws.Cells.Item(98, 14).Value = -14418.5
$ws.Cells.Item(116, 8).Value = 2967.3
$ws.Cells.Item(116, 9).Value = 2531
$ws.Cells.Item(116, 11).Value = 7593
$ws.Cells.Item(116, 13).Value = -4151
$ws.Cells.Item(122, 8).Value = 726.5417
$ws.Cells.Item(122, 9).Value = 726.3333
$ws.Cells.Item(122, 10).Value = 726.75
$ws.Cells.Item(122, 11).Value = 6536.9997
$ws.Cells.Item(122, 12).Value = 6540.75
$ws.Cells.Item(122, 13).Value = -4086.9997
$ws.Cells.Item(122, 14).Value = -11440.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4430.619
$ws.Cells.Item(80, 10).Value = 4806
$ws.Cells.Item(80, 12).Value = 4806
$ws.Cells.Item(80, 14).Value = -6802
$ws.Cells.Item(83, 8).Value = 4430.619
$ws.Cells.Item(83, 10).Value = 4806
$ws.Cells.Item(83, 12).Value = 24030
$ws.Cells.Item(83, 14).Value = -34014
$ws.Cells.Item(126, 8).Value = 7003.8335
$ws.Cells.Item(126, 9).Value = 6505.75
$ws.Cells.Item(126, 11).Value = 19517.25
$ws.Cells.Item(126, 13).Value = -17047.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 5962.6
$ws.Cells.Item(61, 9).Value = 4391.1875
$ws.Cells.Item(61, 10).Value = 12248.25
$ws.Cells.Item(61, 11).Value = 4391.1875
$ws.Cells.Item(61, 12).Value = 12248.25
$ws.Cells.Item(61, 13).Value = -4189.1875
$ws.Cells.Item(61, 14).Value = -12652.25
$ws.Cells.Item(82, 8).Value = 7426.5
$ws.Cells.Item(82, 9).Value = 9244.333000000001
$ws.Cells.Item(82, 10).Value = 1973
$ws.Cells.Item(82, 11).Value = 9244.333000000001
$ws.Cells.Item(82, 12).Value = 1973
$ws.Cells.Item(82, 13).Value = -8883.333000000001
$ws.Cells.Item(82, 14).Value = -2695
$ws.Cells.Item(85, 8).Value = 7426.5
$ws.Cells.Item(85, 9).Value = 9244.333000000001
$ws.Cells.Item(85, 10).Value = 1973
$ws.Cells.Item(85, 11).Value = 9244.333000000001
$ws.Cells.Item(85, 12).Value = 1973
$ws.Cells.Item(85, 13).Value = -7996.333000000001
$ws.Cells.Item(85, 14).Value = -4469
$ws.Cells.Item(113, 8).Value = 5962.6
$ws.Cells.Item(113, 9).Value = 4391.1875
$ws.Cells.Item(113, 10).Value = 12248.25
$ws.Cells.Item(113, 11).Value = 4391.1875
$ws.Cells.Item(113, 12).Value = 12248.25
$ws.Cells.Item(113, 13).Value = -2221.1875
$ws.Cells.Item(113, 14).Value = -16588.25
$ws.Cells.Item(122, 8).Value = 6101.6875
$ws.Cells.Item(122, 9).Value = 5794.615
$ws.Cells.Item(122, 11).Value = 17383.845
$ws.Cells.Item(122, 13).Value = -14933.845
